# Update "想去人数" (interested-count) figures in the "展览" and "全部类型"
# sheets to reflect newly generated site data (gh-pages output refresh).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new F-column value
$updates = @{
    2  = 2798
    4  = 92
    5  = 6674
    6  = 1549
    7  = 17
    9  = 36
    10 = 91
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
